{"js": "// Update the worksheet date and the 25 multiplication prompts in the\n// generated practice sheet (output regenerated at c8c62b6).\nconst replacements = [\n  [\"2025-10-12 Sunday\", \"2025-10-13 Monday\"],\n  [\"950\u00d77=\", \"845\u00d77=\"],\n  [\"647\u00d73=\", \"145\u00d72=\"],\n  [\"862\u00d75=\", \"984\u00d74=\"],\n  [\"900\u00d76=\", \"839\u00d73=\"],\n  [\"417\u00d75=\", \"551\u00d79=\"],\n  [\"975\u00d74=\", \"790\u00d73=\"],\n  [\"469\u00d76=\", \"523\u00d75=\"],\n  [\"225\u00d78=\", \"645\u00d77=\"],\n  [\"889\u00d78=\", \"337\u00d79=\"],\n  [\"193\u00d76=\", \"888\u00d78=\"],\n  [\"407\u00d77=\", \"719\u00d73=\"],\n  [\"666\u00d75=\", \"593\u00d76=\"],\n  [\"458\u00d78=\", \"828\u00d72=\"],\n  [\"449\u00d77=\", \"810\u00d79=\"],\n  [\"641\u00d74=\", \"891\u00d78=\"],\n  [\"775\u00d78=\", \"930\u00d76=\"],\n  [\"603\u00d72=\", \"822\u00d72=\"],\n  [\"285\u00d79=\", \"656\u00d74=\"],\n  [\"177\u00d79=\", \"885\u00d79=\"],\n  [\"848\u00d72=\", \"159\u00d79=\"],\n  [\"850\u00d79=\", \"749\u00d75=\"],\n  [\"212\u00d72=\", \"198\u00d75=\"],\n  [\"212\u00d75=\", \"506\u00d79=\"],\n  [\"326\u00d79=\", \"982\u00d72=\"],\n  [\"441\u00d76=\", \"897\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 multiplication prompts in the\n# generated practice sheet (output regenerated at c8c62b6).\n$d = $word.ActiveDocument\n\n# old-text -> new-text pairs, exactly as they appear in the document body.\n$replacements = @(\n    @(\"2025-10-12 Sunday\", \"2025-10-13 Monday\"),\n    @(\"950\u00d77=\", \"845\u00d77=\"),\n    @(\"647\u00d73=\", \"145\u00d72=\"),\n    @(\"862\u00d75=\", \"984\u00d74=\"),\n    @(\"900\u00d76=\", \"839\u00d73=\"),\n    @(\"417\u00d75=\", \"551\u00d79=\"),\n    @(\"975\u00d74=\", \"790\u00d73=\"),\n    @(\"469\u00d76=\", \"523\u00d75=\"),\n    @(\"225\u00d78=\", \"645\u00d77=\"),\n    @(\"889\u00d78=\", \"337\u00d79=\"),\n    @(\"193\u00d76=\", \"888\u00d78=\"),\n    @(\"407\u00d77=\", \"719\u00d73=\"),\n    @(\"666\u00d75=\", \"593\u00d76=\"),\n    @(\"458\u00d78=\", \"828\u00d72=\"),\n    @(\"449\u00d77=\", \"810\u00d79=\"),\n    @(\"641\u00d74=\", \"891\u00d78=\"),\n    @(\"775\u00d78=\", \"930\u00d76=\"),\n    @(\"603\u00d72=\", \"822\u00d72=\"),\n    @(\"285\u00d79=\", \"656\u00d74=\"),\n    @(\"177\u00d79=\", \"885\u00d79=\"),\n    @(\"848\u00d72=\", \"159\u00d79=\"),\n    @(\"850\u00d79=\", \"749\u00d75=\"),\n    @(\"212\u00d72=\", \"198\u00d75=\"),\n    @(\"212\u00d75=\", \"506\u00d79=\"),\n    @(\"326\u00d79=\", \"982\u00d72=\"),\n    @(\"441\u00d76=\", \"897\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
